$d = $word.ActiveDocument

# --- Intercept row ---
$d.Content.Find.Execute("400264.64 (78.33-2045290078.92)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "39172449945.7 (14211.69-107973108364167344)", 2)
$d.Content.Find.Execute("0.0030669466", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0012653539", 2)

# --- Mean age row ---
$d.Content.Find.Execute("0.97 (0.9-1.04)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.6 (0.43-0.83)", 2)
$d.Content.Find.Execute("0.3991497886", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0021128004", 2)

# --- Female sex percentage row ---
$d.Content.Find.Execute("0.91 (0.84-0.99)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.97 (0.93-1.01)", 2)
$d.Content.Find.Execute("0.0264039331", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.2015492312", 2)

# --- Ethnicity white percentage row ---
$d.Content.Find.Execute("0.95 (0.91-0.99)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.03 (1.01-1.05)", 2)
$d.Content.Find.Execute("0.0152487680", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0004991693", 2)

# --- Mean CPD row ---
$d.Content.Find.Execute("0.85 (0.77-0.93)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.91 (0.84-1)", 2)
$d.Content.Find.Execute("0.0007705366", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0507243308", 2)

# --- Behavioural support only row (label + values) ---
$d.Content.Find.Execute("Behavioural support only", $true, $false, $false, $false, $false, `
    $true, 1, $false, "EMA study type - Interventional", 2)
$d.Content.Find.Execute("101 (4.41-2314.94)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.84 (0.96-3.56)", 2)
$d.Content.Find.Execute("0.0038755014", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0681411829", 2)

# --- Remove "Pharmacological support only" row entirely ---
$t = $d.Tables.Item(1)
$t.Rows.Item(8).Delete()

# --- Study duration days row (pval only changes) ---
$d.Content.Find.Execute("0.2288382263", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.2288382384", 2)

# --- Quality 3 not reported row (label + values) ---
$d.Content.Find.Execute("Quality 3 not reported", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Incentive schedule - Flat payment", 2)
$d.Content.Find.Execute("0.86 (0.32-2.31)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0 (0-0.04)", 2)
$d.Content.Find.Execute("0.7585498702", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0019732057", 2)
